# Apply edits described by the diff for data/case1/20/Q1_1.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A width (15.42578125 -> 16.42578125)
$ws.Columns.Item(1).ColumnWidth = 16.42578125

# Update cell values in column A (rows 1-33)
$ws.Cells.Item(1, 1).Value = 0.04517650658601724
$ws.Cells.Item(2, 1).Value = -0.0059999999400304205
$ws.Cells.Item(3, 1).Value = 0.039659189933367145
$ws.Cells.Item(4, 1).Value = -0.0079999999189439563
$ws.Cells.Item(5, 1).Value = -0.0029999999688801182
$ws.Cells.Item(6, 1).Value = -0.0019999999810274005
$ws.Cells.Item(7, 1).Value = -0.0099999999025199848
$ws.Cells.Item(8, 1).Value = -0.0099999999022961639
$ws.Cells.Item(9, 1).Value = -0.0019999999836399773
$ws.Cells.Item(10, 1).Value = 0.014961780460239993
$ws.Cells.Item(11, 1).Value = 0.015181246545444615
$ws.Cells.Item(12, 1).Value = -0.0034999999698013262
$ws.Cells.Item(13, 1).Value = -0.0034999999655846992
$ws.Cells.Item(14, 1).Value = -0.0079999999201909588
$ws.Cells.Item(15, 1).Value = -0.00099999998922850608
$ws.Cells.Item(16, 1).Value = -0.0019999999787141398
$ws.Cells.Item(17, 1).Value = -0.0019999999781186162
$ws.Cells.Item(18, 1).Value = -0.0039999999579771739
$ws.Cells.Item(19, 1).Value = -0.0039999999600270897
$ws.Cells.Item(20, 1).Value = -0.0039999999567932321
$ws.Cells.Item(21, 1).Value = -0.0039999999562656541
$ws.Cells.Item(22, 1).Value = -0.0039999999558339994
$ws.Cells.Item(23, 1).Value = -0.0049999999453680388
$ws.Cells.Item(24, 1).Value = -0.019999999789913403
$ws.Cells.Item(25, 1).Value = -0.019999999786656453
$ws.Cells.Item(26, 1).Value = 0.074687587263193578
$ws.Cells.Item(27, 1).Value = -0.0024999999743009482
$ws.Cells.Item(28, 1).Value = -0.0019999999688931069
$ws.Cells.Item(29, 1).Value = -0.0069999999124323864
$ws.Cells.Item(30, 1).Value = -0.059999999382528646
$ws.Cells.Item(31, 1).Value = -0.0069999999113399269
$ws.Cells.Item(32, 1).Value = 0.024728083403855905
$ws.Cells.Item(33, 1).Value = -0.0039999999422253296
